$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("avatar_items")

# Insert a new column H for "image_url" (old H/I/J -> I/J/K)
$ws.Columns.Item(8).Insert()

$ws.Range("H5").Value = "TEXT"
$ws.Range("H6").Value = "image_url"
$ws.Hyperlinks.Add($ws.Range("H7"), "https://i.rtings.com/assets/products/NNCSyYNT/keychron-c1/design-medium.jpg")
$ws.Range("H8").Value = $false

$ws.Columns.Item(8).AutoFit()
$ws.PageSetup.Orientation = 1

$ws.Activate()
$ws.Range("H17").Select() | Out-Null
